$d = $word.ActiveDocument
$d.Content.Find.Execute("Computing-2022.git", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Computing-2023.git", 2)
